# Update monthly "worked hours" figures and restore each sheet's
# selection/active-cell state.

$wb = $excel.ActiveWorkbook

# --- 202301 -----------------------------------------------------------
$ws1 = $wb.Worksheets.Item("202301")
$ws1.Range("B10").Value = 10   # was 1
$ws1.Range("B14").Value = 0    # was 9
$ws1.Range("B15").Value = 10   # was 0
$ws1.Range("B25").Select()

# --- 202302 -------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("202302")
$ws2.Range("A2").Select()

# --- 202304 -------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("202304")
$ws4.Range("C13").Select()

# --- 202303 (active tab - select this one last so it stays active) ------
$ws3 = $wb.Worksheets.Item("202303")
$ws3.Range("B22").Value = 12   # was 3
$ws3.Range("B23").Select()
